$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 14.89044
$ws.Range("H2").Value = 44.67131999999999
$ws.Range("I2").Value = 0.2565914865678757
$ws.Range("J2").Value = 0.2565914865678757
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 908.93637661012
$ws.Range("R2").Value = 8180.42738949108
$ws.Range("S2").Value = 0.05243738158189405
$ws.Range("T2").Value = 0.05243738158189405
$ws.Range("G3").Value = 14.89044
$ws.Range("H3").Value = 44.67131999999999
$ws.Range("I3").Value = 0.2565914865678757
$ws.Range("J3").Value = 0.2565914865678757
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 1583.06917710504
$ws.Range("R3").Value = 14247.62259394536
$ws.Range("S3").Value = 0.09132872734171495
$ws.Range("T3").Value = 0.09132872734171496
$ws.Range("G4").Value = 14.89044
$ws.Range("H4").Value = 44.67131999999999
$ws.Range("I4").Value = 0.2565914865678757
$ws.Range("J4").Value = 0.2565914865678757
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 1955.686703873439
$ws.Range("R4").Value = 17601.18033486096
$ws.Range("S4").Value = 0.1128253776442667
$ws.Range("T4").Value = 0.1128253776442667
$ws.Range("I5").Value = 0.6332641083323323
$ws.Range("J5").Value = 0.6332641083323323
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 2243.241939800552
$ws.Range("R5").Value = 20189.17745820497
$ws.Range("S5").Value = 0.1294147055886684
$ws.Range("T5").Value = 0.1294147055886684
$ws.Range("I6").Value = 0.6332641083323323
$ws.Range("J6").Value = 0.6332641083323323
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.2253979890711564
$ws.Range("T6").Value = 0.2253979890711564
$ws.Range("I7").Value = 0.6332641083323323
$ws.Range("J7").Value = 0.6332641083323323
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 4826.606733026594
$ws.Range("R7").Value = 43439.46059723935
$ws.Range("S7").Value = 0.2784514136725074
$ws.Range("T7").Value = 0.2784514136725075
$ws.Range("G8").Value = 6.391867
$ws.Range("H8").Value = 19.175601
$ws.Range("I8").Value = 0.1101444050997921
$ws.Range("J8").Value = 0.1101444050997921
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 390.1698291490244
$ws.Range("R8").Value = 3511.528462341219
$ws.Range("S8").Value = 0.02250925888689095
$ws.Range("T8").Value = 0.02250925888689095
$ws.Range("G9").Value = 6.391867
$ws.Range("H9").Value = 19.175601
$ws.Range("I9").Value = 0.1101444050997921
$ws.Range("J9").Value = 0.1101444050997921
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 679.547926848022
$ws.Range("R9").Value = 6115.931341632198
$ws.Range("S9").Value = 0.03920374941556499
$ws.Range("T9").Value = 0.039203749415565
$ws.Range("G10").Value = 6.391867
$ws.Range("H10").Value = 19.175601
$ws.Range("I10").Value = 0.1101444050997921
$ws.Range("J10").Value = 0.1101444050997921
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 839.4976444502253
$ws.Range("R10").Value = 7555.478800052028
$ws.Range("S10").Value = 0.04843139679733614
$ws.Range("T10").Value = 0.04843139679733614
